$wb = $excel.ActiveWorkbook

# Add the new "Translations question" worksheet right after the existing
# "Translations" sheet (this also makes it the active/selected sheet, as in
# the target workbook where bookViews/workbookView gets activeTab="1").
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item(1))
$ws.Name = "Translations question"

# Header row
$ws.Range("A1").Value = "Entity Id"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Index"
$ws.Range("D1").Value = "Original"
$ws.Range("E1").Value = "Translation"

# Data row - write the text columns first so the new shared strings are
# interned in the same order as the target file (Combobox Option, then the
# Russian translation, then the long numeric-looking entity id).
$ws.Range("B2").Value = "OptionTitle"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Combobox Option"
$ws.Range("E2").Value = "Опция Комбобокса"

# A2 holds a 34-digit number-looking identifier that must be stored as text
# (format code 49, "@") so Excel doesn't coerce/round it as a number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "11111111111111111111111111111111"

# Column widths matching the target layout.
$ws.Columns.Item(1).ColumnWidth = 42.42
$ws.Columns.Item(2).ColumnWidth = 10.26
$ws.Columns.Item(3).ColumnWidth = 10.6
$ws.Columns.Item(4).ColumnWidth = 15.92
$ws.Columns.Item(5).ColumnWidth = 17.6

# Match page setup / selection state of the new sheet.
$ws.PageSetup.Orientation = 1
$ws.Range("A3").Select() | Out-Null
